# Update the two keyword cells on the "addMultipleCustomerTest" sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("addMultipleCustomerTest")
$ws1.Range("A4").Value = "No"
$ws1.Range("B4").Value = "Cust"

# Move the selection on that sheet down to A5 (no longer the active tab)
[void]$ws1.Range("A5").Select()

# Switch the active sheet to "openAccountMethod" (its own remembered
# selection, A4, is left untouched)
$ws2 = $wb.Worksheets.Item("openAccountMethod")
$ws2.Activate()
